$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6:D6").Value = "Enter New Application"
$ws.Range("A7:D7").Value = "Reconcile Applications"
